# Apply cryptos.xlsx price/volume/coin updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (quote-prefixed) so Excel does not
# auto-convert strings like '1.001' or '290.00' into numbers, then strip the
# quote-prefix formatting again so the cell style is left untouched.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '22.036.93'
Set-TextValue $ws.Range("E2") '  -1.00%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.553.47'
Set-TextValue $ws.Range("E3") '  -0.29%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.001'
Set-TextValue $ws.Range("E4") '  +0.07%  '

# Row 5
Set-TextValue $ws.Range("E5") '  -0.04%  '

# Row 6
Set-TextValue $ws.Range("D6") '290.00'
Set-TextValue $ws.Range("E6") '  +0.45%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.3941'
Set-TextValue $ws.Range("E7") '  +3.75%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.3218'
Set-TextValue $ws.Range("E8") '  -2.43%  '

# Row 9
Set-TextValue $ws.Range("D9") '43.67'
Set-TextValue $ws.Range("E9") '  -1.99%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.07232'
Set-TextValue $ws.Range("E10") '  -2.08%  '

# Row 11
Set-TextValue $ws.Range("D11") '1.071'
Set-TextValue $ws.Range("E11") '  -5.97%  '

# Row 12
Set-TextValue $ws.Range("E12") '  +0.04%  '

# Row 13
Set-TextValue $ws.Range("D13") '5.666'
Set-TextValue $ws.Range("E13") '  -3.16%  '

# Row 14
Set-TextValue $ws.Range("D14") '18.68'
Set-TextValue $ws.Range("E14") '  -7.73%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.00001131'
Set-TextValue $ws.Range("E15") '  +5.18%  '

# Row 16
Set-TextValue $ws.Range("B16") 'Chainlink'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D16") '6.617'
Set-TextValue $ws.Range("E16") '  -2.10%  '

# Row 17
Set-TextValue $ws.Range("B17") 'WrappedEther'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D17") '1.551.82'
Set-TextValue $ws.Range("E17") '  -0.58%  '

# Row 18
Set-TextValue $ws.Range("D18") '0.06573'
Set-TextValue $ws.Range("E18") '  -1.34%  '

# Row 19
Set-TextValue $ws.Range("D19") '83.30'
Set-TextValue $ws.Range("E19") '  -3.75%  '

# Row 20
Set-TextValue $ws.Range("D20") '1.000'
Set-TextValue $ws.Range("E20") '  -0.16%  '

# Row 21
Set-TextValue $ws.Range("D21") '6.266'
Set-TextValue $ws.Range("E21") '  -2.34%  '

# Row 22
Set-TextValue $ws.Range("D22") '15.42'
Set-TextValue $ws.Range("E22") '  -4.63%  '

# Row 23
Set-TextValue $ws.Range("D23") '11.28'
Set-TextValue $ws.Range("E23") '  -3.91%  '

# Row 24
Set-TextValue $ws.Range("D24") '22.051.10'
Set-TextValue $ws.Range("E24") '  -0.91%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.359'
Set-TextValue $ws.Range("E25") '  +3.53%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.418'
Set-TextValue $ws.Range("E26") '  -5.74%  '

# Row 27
Set-TextValue $ws.Range("D27") '148.81'
Set-TextValue $ws.Range("E27") '  -1.37%  '

# Row 28
Set-TextValue $ws.Range("D28") '18.49'
Set-TextValue $ws.Range("E28") '  -4.21%  '

# Row 29
Set-TextValue $ws.Range("D29") '4.874'
Set-TextValue $ws.Range("E29") '  -1.29%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.725.65'
Set-TextValue $ws.Range("E30") '  -0.12%  '

# Row 31
Set-TextValue $ws.Range("D31") '118.37'
Set-TextValue $ws.Range("E31") '  -3.62%  '

# Row 32
Set-TextValue $ws.Range("D32") '0.9793'
Set-TextValue $ws.Range("E32") '  -10.22%  '

# Row 33
Set-TextValue $ws.Range("D33") '5.795'
Set-TextValue $ws.Range("E33") '  -2.11%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.08297'
Set-TextValue $ws.Range("E34") '  +0.98%  '

# Row 35
Set-TextValue $ws.Range("D35") '1.604'
Set-TextValue $ws.Range("E35") '  -16.58%  '

# Row 36
Set-TextValue $ws.Range("D36") '9.018'
Set-TextValue $ws.Range("E36") '  -3.91%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.02255'
Set-TextValue $ws.Range("E37") '  -4.10%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.06040'
Set-TextValue $ws.Range("E38") '  -4.37%  '

# Row 39
Set-TextValue $ws.Range("D39") '5.081'
Set-TextValue $ws.Range("E39") '  -5.04%  '

# Row 40
Set-TextValue $ws.Range("D40") '1.207'
Set-TextValue $ws.Range("E40") '  -2.05%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.2034'
Set-TextValue $ws.Range("E41") '  -5.65%  '

# Row 42
Set-TextValue $ws.Range("E42") '  -0.07%  '

# Row 43
Set-TextValue $ws.Range("E43") '  -3.59%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.5788'
Set-TextValue $ws.Range("E44") '  -4.74%  '

# Row 45
Set-TextValue $ws.Range("D45") '3.742'
Set-TextValue $ws.Range("E45") '  -0.19%  '

# Row 46
Set-TextValue $ws.Range("D46") '12.98'
Set-TextValue $ws.Range("E46") '  -5.49%  '

# Row 47
Set-TextValue $ws.Range("D47") '0.5562'
Set-TextValue $ws.Range("E47") '  -5.59%  '

# Row 48
Set-TextValue $ws.Range("B48") 'Quant'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D48") '117.49'
Set-TextValue $ws.Range("E48") '  -4.51%  '

# Row 49
Set-TextValue $ws.Range("B49") 'NEARProtocol'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D49") '1.888'
Set-TextValue $ws.Range("E49") '  -4.21%  '

# Row 50
Set-TextValue $ws.Range("D50") '1.131'
Set-TextValue $ws.Range("E50") '  -3.95%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.06820'
Set-TextValue $ws.Range("E51") '  -3.47%  '
